$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "sd_PAR"

$values = @(
    0.158199379687724,
    281.299233161912,
    256.145614053599,
    177.929878766903,
    0.0124225998749988,
    262.998573047922,
    251.952495046615,
    84.8061081411872,
    0.0124225998749988,
    249.293041647855,
    215.085800210739,
    150.35775651543,
    1.21174758809033,
    244.251964568018,
    249.213167688364,
    157.347737696291
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 6).Value = $v
    $row = $row + 1
}
